$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New columns D (fecha_inicio) and E (fecha_fin) ------------------------
# Give the two new columns their final width/number-format before writing
# headers into them, matching how the header columns A:C already carry a
# text-style width.
$ws.Columns("D:E").ColumnWidth = 10.721354166666666

$ws.Range("D1").Value = "fecha_inicio"
$ws.Range("E1").Value = "fecha_fin"

# Header cells: bold (inherited from the row style) + text number format,
# same treatment the rest of the header row already has.
$ws.Range("D1:E1").Font.Bold = $true
$ws.Range("D1:E1").NumberFormat = "@"

# Whole columns get the text number format too, so new rows typed under the
# headers keep dates as plain YYYY-MM-DD text instead of being reinterpreted.
$ws.Columns("D:E").NumberFormat = "@"

# --- Comments describing the expected date format --------------------------
$commentText = "UTEG:" + [char]10 + "Formato ==> YYYY-MM-DD"
$null = $ws.Range("D1").AddComment($commentText)
$null = $ws.Range("E1").AddComment($commentText)

# --- Selection matches the post-edit saved state ----------------------------
$null = $ws.Range("E19").Select()
